$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Estabilizacion de pagos": update Ahorros account numbers so this
# payments datadriven sheet uses its own unique numeroCuenta values
# (previously duplicated with another sheet) instead of colliding ones.
$ws.Range("T2").Value = "406-125170-01"
$ws.Range("T4").Value = "406-125170-01"
$ws.Range("T3").Value = "406-725170-07"
$ws.Range("T5").Value = "406-725170-07"

# Leave the cursor/selection where the author last left it when saving.
[void]$ws.Range("R10").Select()
